$wb = $excel.ActiveWorkbook

# "Code presentation" sheet: move Songtao from week 3 (B4) to week 9 (B10)
$wsCode = $wb.Worksheets.Item("Code presentation")
$wsCode.Range("B4").ClearContents()
$wsCode.Range("B10").Value = "Songtao"
$wsCode.Range("B10").Select()

# "Lead discussion" sheet: move Songtao from week 11 (B12) to week 5 (B6)
$wsLead = $wb.Worksheets.Item("Lead discussion")
$wsLead.Range("B12").ClearContents()
$wsLead.Range("B6").Value = "Songtao"
$wsLead.Range("B6").Select()
